# Add the "2022-Q4" data and renumber/relabel the existing quarter rollup.
#
# Final sheet tab order: 总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet and move it into position (right
#    after "总计", i.e. right before the existing "2022-Q3" sheet).
#    NOTE: worksheet object references captured before Add()/Move() go stale,
#    so every sheet handle used below is (re-)fetched *after* the sheet list
#    has been finalised.
# ---------------------------------------------------------------------------
$newSheet = $sheets.Add()
$newSheet.Name = "2022-Q4"
$beforeThis = $sheets.Item("2022-Q3")
$newSheet.Move($beforeThis)

$q4 = $wb.Worksheets.Item("2022-Q4")
$q3 = $wb.Worksheets.Item("2022-Q3")
$total = $wb.Worksheets.Item("总计")

function Set-TextCell($ws, $addr, $val) {
    # Force text storage so numeric-looking strings (fund codes with leading
    # zeros, "9.60"-style decimals, ...) keep their literal representation
    # instead of being coerced to a Double by the normal Value setter.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# ---------------------------------------------------------------------------
# 2. Populate "2022-Q4" — same layout as the other quarter sheets.
#    Pull the header-row and index-column formatting (bold/border style)
#    from "2022-Q3" so the new sheet matches the existing look.
# ---------------------------------------------------------------------------
$q3.Range("A1:H1").Copy($q4.Range("A1:H1"))
$q3.Range("A2").Copy($q4.Range("A2:A9"))

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
Set-TextCell $q4 "B2" "016250"
$q4.Range("C2").Value = "华夏远见成长一年持有混合A"
Set-TextCell $q4 "D2" "9.60"
Set-TextCell $q4 "E2" "88.62"
Set-TextCell $q4 "F2" "4.66"
Set-TextCell $q4 "G2" "0.4474"
$q4.Range("H2").Value = 4

$q4.Range("A3").Value = 1
Set-TextCell $q4 "B3" "003501"
$q4.Range("C3").Value = "泰达宏利睿智稳健灵活配置混合A"
Set-TextCell $q4 "D3" "9.84"
Set-TextCell $q4 "E3" "82.46"
Set-TextCell $q4 "F3" "2.10"
Set-TextCell $q4 "G3" "0.2066"
$q4.Range("H3").Value = 9

$q4.Range("A4").Value = 2
Set-TextCell $q4 "B4" "016251"
$q4.Range("C4").Value = "华夏远见成长一年持有混合C"
Set-TextCell $q4 "D4" "2.97"
Set-TextCell $q4 "E4" "88.62"
Set-TextCell $q4 "F4" "4.66"
Set-TextCell $q4 "G4" "0.1384"
$q4.Range("H4").Value = 4

$q4.Range("A5").Value = 3
Set-TextCell $q4 "B5" "013280"
$q4.Range("C5").Value = "泰达宏利睿智稳健灵活配置混合C"
Set-TextCell $q4 "D5" "5.43"
Set-TextCell $q4 "E5" "82.46"
Set-TextCell $q4 "F5" "2.10"
Set-TextCell $q4 "G5" "0.1140"
$q4.Range("H5").Value = 9

$q4.Range("A6").Value = 4
Set-TextCell $q4 "B6" "011431"
$q4.Range("C6").Value = "泰达宏利消费服务混合A"
Set-TextCell $q4 "D6" "2.00"
Set-TextCell $q4 "E6" "87.31"
Set-TextCell $q4 "F6" "3.30"
Set-TextCell $q4 "G6" "0.0660"
$q4.Range("H6").Value = 7

$q4.Range("A7").Value = 5
Set-TextCell $q4 "B7" "011432"
$q4.Range("C7").Value = "泰达宏利消费服务混合C"
Set-TextCell $q4 "D7" "1.39"
Set-TextCell $q4 "E7" "87.31"
Set-TextCell $q4 "F7" "3.30"
Set-TextCell $q4 "G7" "0.0459"
$q4.Range("H7").Value = 7

$q4.Range("A8").Value = 6
Set-TextCell $q4 "B8" "001744"
$q4.Range("C8").Value = "诺安进取回报灵活配置混合"
Set-TextCell $q4 "D8" "0.59"
Set-TextCell $q4 "E8" "69.55"
Set-TextCell $q4 "F8" "3.48"
Set-TextCell $q4 "G8" "0.0205"
$q4.Range("H8").Value = 7

$q4.Range("A9").Value = 7
Set-TextCell $q4 "B9" "519139"
$q4.Range("C9").Value = "海富通沪港深灵活配置混合"
Set-TextCell $q4 "D9" "0.67"
Set-TextCell $q4 "E9" "92.35"
Set-TextCell $q4 "F9" "2.93"
Set-TextCell $q4 "G9" "0.0196"
$q4.Range("H9").Value = 7

# ---------------------------------------------------------------------------
# 3. Update "总计": add a row for 2022-Q4 at the top of the data and shift
#    the other quarters' figures down one row, appending 2021-Q4 at the end.
# ---------------------------------------------------------------------------
$total.Range("A5").Copy($total.Range("A6"))

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 8
$total.Range("D2").Value = 1.06

$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 20
$total.Range("D3").Value = 1.56

$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 0.62

$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 6
$total.Range("D5").Value = 0.07000000000000001

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q4"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.17
